$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values that changed, preserving them as text
# by setting NumberFormat to Text ("@") before assigning the value.
$priceUpdates = @{
    "D2" = "243.72"
    "D3" = "22.97"
    "D4" = "5.415"
    "D5" = "0.05951"
    "D6" = "3.452"
    "D7" = "6.525"
    "D8" = "0.8128"
    "D9" = "0.9178"
    "D10" = "0.1413"
    "D11" = "0.07492"
    "D12" = "0.03277"
    "D13" = "0.03058"
    "D14" = "0.09348"
    "D15" = "3.845"
    "D16" = "0.001558"
    "D17" = "0.04666"
    "D18" = "0.0005938"
    "D19" = "0.006121"
    "D20" = "0.004996"
    "D21" = "0.0009799"
    "D22" = "0.00007896"
    "D23" = "3.609"
    "D26" = "0.1301"
    "D27" = "0.0002393"
    "D40" = "0.03930"
    "D41" = "0.006162"
    "D43" = "0.002999"
    "D44" = "0.008544"
    "D45" = "0.00005235"
    "D48" = "0.8997"
    "D49" = "0.002281"
    "D50" = "0.00002099"
    "D51" = "0.0001999"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Update Hora column (G) for rows 2-51 from "5" to "6", preserving as text
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = "6"
}
